$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.802.84"
$ws.Range("E2").Value = "  -0.66%  "
$ws.Range("D3").Value = "2.749.04"
$ws.Range("E3").Value = "  +0.24%  "
$ws.Range("E4").Value = "  -0.02%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "572.16"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -1.74%  "
$ws.Range("E6").Value = "  -1.24%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.997"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("E9").Value = "  -4.27%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.161"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +1.09%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.381"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -3.67%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "5.54"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -18.15%  "
$ws.Range("D13").Value = "3.236.16"
$ws.Range("E13").Value = "  +0.11%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "26.37"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -3.53%  "
$ws.Range("D15").Value = "63.446.76"
$ws.Range("E15").Value = "  -1.09%  "
$ws.Range("E16").Value = "  -3.86%  "
$ws.Range("D17").Value = "2.752.46"
$ws.Range("E17").Value = "  -0.44%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "12.07"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -0.35%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "4.79"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -3.58%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "354.06"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -2.71%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "6.71"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -5.15%  "
$ws.Range("E22").Value = "  +0.01%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "0.533"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -0.45%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "64.94"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -3.69%  "
$ws.Range("E25").Value = "  -1.49%  "
$ws.Range("E26").Value = "  +0.37%  "
$ws.Range("E27").Value = "  -2.98%  "
$ws.Range("D28").Value = "0.0₃0898"
$ws.Range("E28").Value = "  -2.62%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "1.93"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -4.83%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "6.98"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -3.42%  "
$ws.Range("B31").Value = "Fetch.AI"
$ws.Range("C31").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "1.21"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -4.22%  "
$ws.Range("B32").Value = "Monero"
$ws.Range("C32").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "168.78"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -3.46%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "20.11"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -3.02%  "
$ws.Range("E34").Value = "  +0.16%  "
$ws.Range("E35").Value = "  -1.64%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "1.44"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -1.45%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "1.79"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -2.25%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.975"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -4.03%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "6.14"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +4.55%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "4.12"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -4.75%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "327.70"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -4.18%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "38.93"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -1.22%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "21.32"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -3.81%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.0583"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -3.44%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "21.28"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -4.16%  "
$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.0253"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -3.03%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "134.81"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -2.79%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "0.623"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -4.82%  "
$ws.Range("E49").Value = "  -1.54%  "
$ws.Range("E50").Value = "  +0.11%  "
$ws.Range("E51").Value = "  +0.25%  "
